$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write new string cells in the exact order the strings were first introduced
# (this preserves the shared-strings table order) before filling in the rest.
$ws.Range("E6").Value = "Fixed controller behaviour. NPC model add, gun functionality"
$ws.Range("B6").Value = "9.15 - 16.00"

$ws.Range("B7").Value = "10.30 - 12.00"
$ws.Range("E7").Value = "Fixed enemy animations"

$ws.Range("B8").Value = "10.15 - 12.30"
$ws.Range("B9").Value = "13.30 - 17.00"

$ws.Range("E8").Value = "Worked on gun system"
$ws.Range("E9").Value = "Worked on gun system"

$ws.Range("B10").Value = "9.00 - 14.00"
$ws.Range("B11").Value = "9.00 - 14.00"
$ws.Range("B12").Value = "9.00 - 14.00"
$ws.Range("B13").Value = "9.00 - 14.00"
$ws.Range("B14").Value = "9.00 - 14.00"
$ws.Range("B15").Value = "9.00 - 14.00"
$ws.Range("B16").Value = "9.00 - 14.00"
$ws.Range("B17").Value = "9.00 - 14.00"
$ws.Range("B18").Value = "9.00 - 14.00"

# Dates (column A), formatted like the existing entries (d-mmm)
$dates = @{
    6  = 44649
    7  = 44651
    8  = 44655
    9  = 44655
    10 = 44656
    11 = 44657
    12 = 44658
    13 = 44659
    14 = 44662
    15 = 44663
    16 = 44664
    17 = 44665
    18 = 44666
}
foreach ($r in $dates.Keys) {
    $cell = $ws.Range("A$r")
    $cell.Value = $dates[$r]
    $cell.NumberFormat = "d-mmm"
}

# Hours (column D)
$hours = @{
    6  = 6.75
    7  = 1.5
    8  = 2.15
    9  = 3.5
    10 = 5
    11 = 5
    12 = 5
    13 = 5
    14 = 5
    15 = 5
    16 = 5
    17 = 5
    18 = 5
}
foreach ($r in $hours.Keys) {
    $ws.Range("D$r").Value = $hours[$r]
}

# Total row
$ws.Range("D30").Formula = "=SUM(D3:D29)"

# Final selection as left by the author
$ws.Range("M29:N30").Select() | Out-Null
